$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a Site id value for row 3, column J (matches pattern of J2 which holds a Site id)
$ws.Range("J3").Value = "a1Zq0000000EJNA"

# Add a new row 4 with an id value in column A (matches pattern of A3 holding an id)
$ws.Range("A4").Value = "001q000000hmj2V"

# Update selection to reflect the new active cell, like the authored workbook
$ws.Range("A4").Select()
